$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Boundary-test fixture update: new VisitsPerMonth ("date") inputs for the
# membership-ranking test rows, plus the recomputed ExpectedRank now that
# the "Gold" tier no longer applies (collapses to "Standard").

# Column B (VisitsPerMonth)
$ws.Range("B2").Value  = 15
$ws.Range("B3").Value  = 15
$ws.Range("B4").Value  = 15
$ws.Range("B5").Value  = 15
$ws.Range("B6").Value  = 15
$ws.Range("B9").Value  = 29
$ws.Range("B10").Value = 30
$ws.Range("B11").Value = 15
$ws.Range("B12").Value = 15
$ws.Range("B13").Value = 15
$ws.Range("B14").Value = 15

# Column D (ExpectedRank): rows that expected "Gold" now expect "Standard"
$ws.Range("D4").Value  = "Standard"
$ws.Range("D5").Value  = "Standard"
$ws.Range("D6").Value  = "Standard"
$ws.Range("D13").Value = "Standard"
$ws.Range("D14").Value = "Standard"

# Reviewer selection now highlights the updated Actual/Result columns
$ws.Range("E2:F14").Select()
